# "Generate Report for Archive"
# - Status text "Ready for handoff" -> "In Translation" everywhere it appears
#   (Overview!E2:F3 and the Status column (C) on the zh-cn / de-de sheets).
# - Narrow the "zh-cn"/"de-de" status columns (Overview E:F, zh-cn/de-de C)
#   from their old width down to the new, narrower width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NB: compare with the literal on the left - PowerShell's -eq coerces
        # the right-hand side to the left-hand side's type, and a boolean
        # cell (e.g. True/False) on the right would otherwise coerce any
        # non-empty string to $true and falsely match.
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value = $newStatus
        }
    }
}

# Overview sheet: columns E (zh-cn) and F (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C (Status)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C (Status)
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
